# Update lamda_1 (B), lamda_2 (C) and the dic_nbre_clients_poisson_2 key/value
# columns (D/E), then trim the now-shorter table (rows 56-58 removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New D (keys) and E (probability values) data for rows 2..55 (54 data rows)
$dVals = @(0,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,50,51,54,55,56,59,62)
$eVals = @(0.136,0.001,0.009000000000000001,0.02,0.026,0.048,0.054,0.057,0.049,0.027,0.033,0.028,0.043,0.027,0.041,0.036,0.028,0.028,0.027,0.021,0.026,0.027,0.012,0.023,0.013,0.012,0.017,0.009000000000000001,0.013,0.012,0.01,0.006,0.007,0.009000000000000001,0.009000000000000001,0.006,0.004,0.002,0.003,0.004,0.005,0.008,0.001,0.004,0.003,0.003,0.003,0.002,0.001,0.001,0.001,0.001,0.002,0.001)

$firstRow = 2
$lastRowBefore = 58
$lastRowAfter = $firstRow + $dVals.Length - 1   # 55

# First, delete the trailing rows that no longer exist in the updated table
# (rows 56 through 58), shrinking the sheet from 58 to 55 used rows.
$ws.Range("A$($lastRowAfter + 1):E$lastRowBefore").EntireRow.Delete() | Out-Null

# Update the auto-correlation scale (lamda_1, column B) and auto capacity
# (lamda_2, column C) values, which are constant across all remaining rows.
$ws.Range("B$firstRow`:B$lastRowAfter").Value = 33.94444444444444
$ws.Range("C$firstRow`:C$lastRowAfter").Value = 1.95

# Update the dictionary keys (D) and probability values (E, "time in ms")
for ($i = 0; $i -lt $dVals.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}
